$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 53
$ws.Cells.Item($row, 1).Value = "2025-04-29 06:42:23"
$ws.Cells.Item($row, 2).Value = 148
